$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2056
$ws1.Range("F4").Value = 855
$ws1.Range("F5").Value = 1157
$ws1.Range("F6").Value = 352

# Sheet "全部类型" - same events appear again, update matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2056
$ws4.Range("F6").Value = 855
$ws4.Range("F7").Value = 1157
$ws4.Range("F8").Value = 352
